# Append district-level Covid positive test data for date serial 44180
# (2020-12-15), matching the author's "Updated: st 16. 12. 2020" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 44180
$startRow = 6150

$districts = @(
    @('Bánovce nad Bebravou', 33),
    @('Banská Bystrica', 63),
    @('Banská Štiavnica', 4),
    @('Bardejov', 48),
    @('Bratislava', 294),
    @('Brezno', 23),
    @('Bytča', 16),
    @('Čadca', 60),
    @('Detva', 11),
    @('Dolný Kubín', 14),
    @('Dunajská Streda', 69),
    @('Galanta', 43),
    @('Gelnica', 7),
    @('Hlohovec', 73),
    @('Humenné', 14),
    @('Ilava', 47),
    @('Kežmarok', 15),
    @('Komárno', 61),
    @('Košice', 98),
    @('Košice - okolie', 38),
    @('Krupina', 1),
    @('Kysucké Nové Mesto', 12),
    @('Levice', 35),
    @('Levoča', 15),
    @('Liptovský Mikuláš', 50),
    @('Lučenec', 92),
    @('Malacky', 52),
    @('Martin', 110),
    @('Medzilaborce', 7),
    @('Michalovce', 26),
    @('Myjava', 31),
    @('Námestovo', 6),
    @('Nitra', 237),
    @('Nové Mesto nad Váhom', 76),
    @('Nové Zámky', 51),
    @('Partizánske', 21),
    @('Pezinok', 32),
    @('Piešťany', 40),
    @('Poltár', 11),
    @('Poprad', 85),
    @('Považská Bystrica', 68),
    @('Prešov', 72),
    @('Prievidza', 98),
    @('Púchov', 52),
    @('Revúca', 16),
    @('Rimavská Sobota', 36),
    @('Rožňava', 31),
    @('Ružomberok', 100),
    @('Sabinov', 6),
    @('Senec', 45),
    @('Senica', 86),
    @('Skalica', 87),
    @('Snina', 19),
    @('Sobrance', 5),
    @('Spišská Nová Ves', 37),
    @('Stará Ľubovňa', 22),
    @('Stropkov', 20),
    @('Svidník', 20),
    @('Šaľa', 35),
    @('Topoľčany', 20),
    @('Trebišov', 77),
    @('Trenčín', 169),
    @('Trnava', 94),
    @('Turčianske Teplice', 11),
    @('Tvrdošín', 11),
    @('Veľký Krtíš', 18),
    @('Vranov nad Topľou', 68),
    @('Zlaté Moravce', 22),
    @('Zvolen', 41),
    @('Žarnovica', 5),
    @('Žiar nad Hronom', 15),
    @('Žilina', 138)
)

$r = $startRow
foreach ($d in $districts) {
    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 2).Value = $d[0]
    $ws.Cells.Item($r, 3).Value = $d[1]
    $r = $r + 1
}

Write-Output "Appended $($districts.Count) rows starting at row $startRow"
